# Add new worksheet "Sheet4" at the end of the workbook (after the current last sheet)
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([Type]::Missing, $lastSheet, 1, [Type]::Missing)
$ws4.Name = "Sheet4"

# Populate the ranking-results table (same layout/shared strings as Sheet1-3)
$ws4.Range("A1").Value = "Model"
$ws4.Range("B1").Value = "Score"
$ws4.Range("A2").Value = "Logistic Regression"
$ws4.Range("B2").Value = 0.763
$ws4.Range("A3").Value = "Neural Network MLP Classification"
$ws4.Range("B3").Value = 0.727
$ws4.Range("A4").Value = "Lasso Regression"
$ws4.Range("B4").Value = 0.706

# Scores are percentages, same formatting as the other ranking sheets
$ws4.Range("B2:B4").NumberFormat = "0.00%"

# Re-sort descending by score (matches sortState on the other sheets)
$ws4.Sort.SortFields.Clear()
$ws4.Sort.SortFields.Add($ws4.Range("B1"), [Type]::Missing, 2, [Type]::Missing, 1) | Out-Null
$ws4.Sort.SetRange($ws4.Range("A1:B4"))
$ws4.Sort.Header = 1
$ws4.Sort.Apply()

# Add the clustered-column "Score" chart, sourced from the new table
$chartObj = $ws4.ChartObjects().Add(60, 100, 420, 225)
$chartObj.Chart.ChartType = 51
$chartObj.Chart.SetSourceData($ws4.Range("A1:B4"))

# Sheet3's selection moves to match the A1:B5 pattern used elsewhere
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1:B5").Select() | Out-Null

# Sheet4 becomes the newly active/selected tab
$ws4.Activate() | Out-Null
$ws4.Range("L16").Select() | Out-Null
